$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.928.99'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '2.617.53'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''308.99'
$ws.Range('E5').Value = '  -1.56%  '
$ws.Range('D6').Value = '''98.71'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('D7').Value = '''0.595'
$ws.Range('E7').Value = '  -1.09%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('D10').Value = '''38.71'
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('D11').Value = '''54.07'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = '''8.04'
$ws.Range('E13').Value = '  -3.79%  '
$ws.Range('D14').Value = '3.014.35'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('D16').Value = '2.621.97'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').Value = '''0.915'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').Value = '''14.81'
$ws.Range('D19').Value = '45.917.52'
$ws.Range('E19').Value = '  -1.37%  '
$ws.Range('E20').Value = '  -1.43%  '
$ws.Range('D21').Value = '''6.74'
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('D22').Value = '''12.69'
$ws.Range('E22').Value = '  -5.43%  '
$ws.Range('D23').Value = '''75.78'
$ws.Range('E23').Value = '  +6.57%  '
$ws.Range('D24').Value = '''282.21'
$ws.Range('E24').Value = '  +10.38%  '
$ws.Range('D25').Value = '''3.02'
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('D26').Value = '''2.24'
$ws.Range('E26').Value = '  +1.10%  '
$ws.Range('D27').Value = '''29.53'
$ws.Range('E27').Value = '  +4.86%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('D30').Value = '''10.46'
$ws.Range('E30').Value = '  -2.32%  '
$ws.Range('D31').Value = '''38.61'
$ws.Range('E31').Value = '  -6.13%  '
$ws.Range('E32').Value = '  -3.71%  '
$ws.Range('D33').Value = '''6.22'
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E34').Value = '  -3.86%  '
$ws.Range('D35').Value = '''2.28'
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').Value = '''156.81'
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.0835'
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').Value = '''2.81'
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('D39').Value = '''0.122'
$ws.Range('E39').Value = '  +2.28%  '
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('D41').Value = '''15.87'
$ws.Range('E41').Value = '  -7.12%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '''22.09'
$ws.Range('E42').Value = '  +4.98%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0326'
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''4.02'
$ws.Range('E44').Value = '  -5.90%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '''3.53'
$ws.Range('E45').Value = '  -3.32%  '
$ws.Range('D46').Value = '2.107.39'
$ws.Range('E46').Value = '  +3.58%  '
$ws.Range('D47').Value = '''0.998'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').Value = '''93.95'
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('D49').Value = '''9.15'
$ws.Range('E49').Value = '  -0.99%  '
$ws.Range('D50').Value = '''109.35'
$ws.Range('E50').Value = '  -3.35%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.865.95'
$ws.Range('E51').Value = '  -0.06%  '
